$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("D4").Value = 44235
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 80
$ws.Range("K4").Value = 14000
$ws.Range("L4").Value = 14000
$ws.Range("M4").Value = 14000
$ws.Range("N4").Value = "$/bandeja 18 kilos"
$ws.Range("P4").Value = 778
$ws.Range("Q4").Value = 18

# Row 5
$ws.Range("D5").Value = 44235
$ws.Range("I5").Value = "Segunda"
$ws.Range("J5").Value = 70
$ws.Range("K5").Value = 12000
$ws.Range("L5").Value = 12000
$ws.Range("M5").Value = 12000
$ws.Range("N5").Value = "$/bandeja 18 kilos"
$ws.Range("P5").Value = 667
$ws.Range("Q5").Value = 18

# Row 6
$ws.Range("D6").Value = 44235
$ws.Range("I6").Value = "Tercera"
$ws.Range("J6").Value = 60
$ws.Range("K6").Value = 10000
$ws.Range("L6").Value = 10000
$ws.Range("M6").Value = 10000
$ws.Range("N6").Value = "$/bandeja 18 kilos"
$ws.Range("P6").Value = 556
$ws.Range("Q6").Value = 18

# Row 7
$ws.Range("D7").Value = 44536
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 87
$ws.Range("K7").Value = 22000
$ws.Range("L7").Value = 22000
$ws.Range("M7").Value = 22000
$ws.Range("N7").Value = "$/bandeja 18 kilos"
$ws.Range("P7").Value = 1222
$ws.Range("Q7").Value = 18

# Row 8
$ws.Range("D8").Value = 44536
$ws.Range("I8").Value = "Segunda"
$ws.Range("J8").Value = 80
$ws.Range("K8").Value = 20000
$ws.Range("L8").Value = 20000
$ws.Range("M8").Value = 20000
$ws.Range("N8").Value = "$/bandeja 18 kilos"
$ws.Range("P8").Value = 1111
$ws.Range("Q8").Value = 18

# Row 9
$ws.Range("D9").Value = 44756
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 65
$ws.Range("K9").Value = 14000
$ws.Range("L9").Value = 14000
$ws.Range("M9").Value = 14000
$ws.Range("N9").Value = "$/caja 15 kilos"
$ws.Range("P9").Value = 933
$ws.Range("Q9").Value = 15

# Row 10
$ws.Range("D10").Value = 44756
$ws.Range("I10").Value = "Segunda"
$ws.Range("J10").Value = 68
$ws.Range("K10").Value = 12000
$ws.Range("L10").Value = 12000
$ws.Range("M10").Value = 12000
$ws.Range("N10").Value = "$/caja 15 kilos"
$ws.Range("P10").Value = 800
$ws.Range("Q10").Value = 15

# Row 11
$ws.Range("D11").Value = 44424
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 75
$ws.Range("K11").Value = 18000
$ws.Range("L11").Value = 18000
$ws.Range("M11").Value = 18000
$ws.Range("N11").Value = "$/caja 15 kilos"
$ws.Range("P11").Value = 1200
$ws.Range("Q11").Value = 15

# Row 12
$ws.Range("D12").Value = 44424
$ws.Range("I12").Value = "Segunda"
$ws.Range("J12").Value = 50
$ws.Range("K12").Value = 12000
$ws.Range("L12").Value = 12000
$ws.Range("M12").Value = 12000
$ws.Range("N12").Value = "$/caja 15 kilos"
$ws.Range("P12").Value = 800
$ws.Range("Q12").Value = 15

# Row 13
$ws.Range("D13").Value = 44242
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 60
$ws.Range("K13").Value = 13000
$ws.Range("L13").Value = 13000
$ws.Range("M13").Value = 13000
$ws.Range("N13").Value = "$/bandeja 18 kilos"
$ws.Range("P13").Value = 722
$ws.Range("Q13").Value = 18

# Row 14
$ws.Range("D14").Value = 44242
$ws.Range("I14").Value = "Segunda"
$ws.Range("J14").Value = 50
$ws.Range("K14").Value = 10000
$ws.Range("L14").Value = 10000
$ws.Range("M14").Value = 10000
$ws.Range("N14").Value = "$/bandeja 18 kilos"
$ws.Range("P14").Value = 556
$ws.Range("Q14").Value = 18

